$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.294.97'
$ws.Range("E2").Value = '  -4.99%  '
$ws.Range("D3").Value = '3.266.90'
$ws.Range("E3").Value = '  -7.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.90'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.07'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = '  -12.51%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.258.64'
$ws.Range("E8").Value = '  -7.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = '  -10.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = '  -14.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.61'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = '  -8.84%  '
$ws.Range("E12").Value = '  -12.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.48'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = '  -16.75%  '
$ws.Range("E14").Value = '  -11.54%  '
$ws.Range("D15").Value = '3.794.86'
$ws.Range("E15").Value = '  -7.45%  '
$ws.Range("D16").Value = '67.304.95'
$ws.Range("E16").Value = '  -5.08%  '
$ws.Range("D17").Value = '3.271.40'
$ws.Range("E17").Value = '  -7.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '535.45'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = '  -11.80%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Value = '  -6.28%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.22'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = '  -14.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.06'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = '  -14.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.760'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = '  -13.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.85'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = '  -13.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.63'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = '  -12.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.58'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = '  -13.51%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -12.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.08'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = '  -11.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.15'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = '  -16.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.33'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = '  -12.90%  '
$ws.Range("E31").Value = '  -11.39%  '
$ws.Range("E32").Value = '  -10.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.59'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = '  -19.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '536.10'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Value = '  -13.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.78'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = '  -15.44%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0458'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null
$ws.Range("E37").Value = '  -8.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.37'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4122) | Out-Null
$ws.Range("E38").Value = '  -6.25%  '
$ws.Range("E39").Value = '  -13.72%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.128'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4122) | Out-Null
$ws.Range("E40").Value = '  -11.03%  '
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.07'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4122) | Out-Null
$ws.Range("E41").Value = '  -16.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").Value = '  -18.90%  '
$ws.Range("D43").Value = '2.929.55'
$ws.Range("E43").Value = '  -12.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.270'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").Value = '  -13.43%  '
$ws.Range("D45").Value = '0.0₃0592'
$ws.Range("E45").Value = '  -18.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.19'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4122) | Out-Null
$ws.Range("E46").Value = '  -12.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.72'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$ws.Range("E47").Value = '  -15.86%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.34'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Value = '  -18.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.36'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4122) | Out-Null
$ws.Range("E50").Value = '  -6.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.113'
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").Value = '  -12.55%  '
$excel.CutCopyMode = $false
